$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "особливый товар"
$ws.Range("A16").Value = "мелочь"
$ws.Range("A17").Value = "деревенский товар"
$ws.Range("A18").Value = "серебреный товар"
$ws.Range("A19").Value = "небогатый товар"
$ws.Range("A20").Value = "крамными товар"
$ws.Range("A21").Value = "мясо"
$ws.Range("A22").Value = "железный товар"
$ws.Range("A25").Value = "щепетильный товар"
$ws.Range("A26").Value = "пушной товар"
$ws.Range("A28").Value = "недорогой товар"
$ws.Range("A29").Value = "питейный припасы"
$ws.Range("A30").Value = "суровский товар"
$ws.Range("A31").Value = "медный товар"
$ws.Range("A32").Value = "внутренний товар"
$ws.Range("A33").Value = "привозный товар"
$ws.Range("A34").Value = "оловянный товар"
$ws.Range("A35").Value = "купецкий товар"
$ws.Range("A36").Value = "галантерейный товар"
$ws.Range("A37").Value = "произрастание"
$ws.Range("A39").Value = "домовый товар"
$ws.Range("A40").Value = "надлежащий товар"
$ws.Range("A41").Value = "рукодельный товар"
$ws.Range("A42").Value = "харчевой припасы"
$ws.Range("A43").Value = "меховой товар"
